$wb = $excel.ActiveWorkbook

# --- RGS565 sheet: update RGB input values (B1:B3) from 55 to 64 ---
$wsRGS = $wb.Worksheets.Item("RGS565")
$wsRGS.Range("B1").Value = 64
$wsRGS.Range("B2").Value = 64
$wsRGS.Range("B3").Value = 64

# Move the selection on RGS565 from D5 to B4
[void]$wsRGS.Range("B4").Select()

# --- APIs sheet becomes the active/selected tab (was gMSS before) ---
$wsAPIs = $wb.Worksheets.Item("APIs")
[void]$wsAPIs.Activate()
